$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# LED2 (row 28) is now DNP (Do Not Populate)
$ws.Range("A28").Copy()
$ws.Range("H28").PasteSpecial(-4122)  # xlPasteFormats - reuse the row's existing cell style
$ws.Range("H28").Value = "DNP"

# R12 (row 44) is now DNP (Do Not Populate)
$ws.Range("A44").Copy()
$ws.Range("H44").PasteSpecial(-4122)  # xlPasteFormats - reuse the row's existing cell style
$ws.Range("H44").Value = "DNP"

$excel.CutCopyMode = $false
